$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,3,4) get their values rotated: row2 -> row3 -> row4 -> row2
# (for columns D, J, K, L, M, O, P). Capture the original values first, then
# write the rotated values back.

# Column indexes: D=4, J=10, K=11, L=12, M=13, O=15, P=16
$colIndexes = @(4, 10, 11, 12, 13, 15, 16)

$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($c in $colIndexes) {
    $orig2[$c] = $ws.Cells.Item(2, $c).Value2
    $orig3[$c] = $ws.Cells.Item(3, $c).Value2
    $orig4[$c] = $ws.Cells.Item(4, $c).Value2
}

foreach ($c in $colIndexes) {
    $ws.Cells.Item(3, $c).Value2 = $orig2[$c]
    $ws.Cells.Item(4, $c).Value2 = $orig3[$c]
    $ws.Cells.Item(2, $c).Value2 = $orig4[$c]
}
